$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string text that changed in place (rows 4 and 6 / Art of Crafting + Luck of the roll) ---
$ws.Range("D4").Value = 'Crafting is the second most set of essential skills behind Enchanting. <br /> <br /> The shop only goes up 2 Billion Gold, while you can at later levels craft to 36 Billion Gold (level 400, top end gear). <br /> <br /> While it is true that monsters drops will scale with you to a max crafting level of 300, you A) have to be in Purgatory to get those kinds of drops and B) it''s all random loot, it could be good, could be bad, could be so-so. You need to spend time leveling crafting. You will out craft the shop some where in the level 200''s of each crafting skill. Each plane has a limit on its max level of gear that is not displayed any where to player.<br /> <br /> For this Quest here’s what I want you to do:<br /> <br /> On PC, either move or teleport to these coordinates: (X/Y) 32/368<br /> <br /> On Mobile, select Map Movement from the primary action drop down. Open Teleport and enter in the coordinates if you do not have the gold. Kill till you do.<br /> <br /> Now that you are here, you should have seen chat update. You got a quest item. Quest items can be seen from the inventory under the quest items tab. This item will give you a slight boost to weapon crafting.<br /> <br /> Now that you have it.<br /> <br /> On Pc, click the blue drop down to select craft. Now select weapons from the list. Select Broken Dagger and Click Craft. Repeat the clicking of the button. Pay attention to the sever messages tab as you craft.<br /> <br /> There is no automation for any form of crafting.<br /> <br /> For Mobile, click the actions drop down and select craft, now follow the same steps above.<br /> You will fail more than you succeed at first, but the quest item should help you a bit. If you would like to see how much:<br /> <br /> Click Character Sheet Tab, Click the crafting skills Tab for the skills section (mobile, select Skill Management, then tap the crafting skills tab) Click on Weapon Crafting to see the skill bonus (applied to your chance to craft) and Skill XP Bonus. You care about the bonus to succeed and the XP bonus to level quickly.<br /> <br /> For now, focus on crafting broken daggers, as you gain levels you unlock new items. Click the help button associated with the crafting action to learn more.<br /> <br /> Because there is not automation for any of the crafting options, it is suggested players put on some Netflix, turn on automation, then go back to craft while automation runs.'
$ws.Range("D6").Value = 'This one is super easy. You have Enchanting to at least level 5 so now what I want you to do is craft yourself a set of gear, a whole set.<br /> <br /> Craft 12 Items and enchant them with: Treasure Goblin Luck (unlocks at level 5 and is a prefix). This gives you 5% towards Skill XP and Bonus for the Skill: Looting. Equip them.<br /> <br /> Your looting should be roughly 60% for Bonus and XP bonus.<br /> <br /> Now change your skill in training to looting at the same 10% XP sacrifice. Explore for 2 hours this time.<br /> <br /> Looting essentially helps you get more drops. In some locations, like Special locations we use your looting skill but capped at 45%. The higher the looting skill the higher the chance for drops and quest drops.<br /> <br /> Enchanting multiple items with the same enchant will stack unless otherwise stated by the enchantment. Caution though, if you have that enchanting help docs page open, look at the intelligence required to enchant. If your intelligence is not high enough, even if you have the level, you wont see the enchantment. Things like Staves and Damage Spells can help raise your int, along with enchants to raise your int (filter the list in the help docs). You can also level your character to raise the int. Heretics are the best enchanters as well as Arcane Alchemists, but every one can do well at enchanting. Theres even enchantments to help with enchanting.'

# --- Numeric tweaks: required_skill_level drops from 10 to 5 for rows 4 and 5 ---
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 5

# --- Column B width: 23 -> 25 characters (raw OOXML width unit) ---
$ws.Columns.Item(2).ColumnWidth = 24.166666666666668

# --- New row 7: "Quests for the wicked" ---
$ws.Range("A7").Value = 21
$ws.Range("B7").Value = 'Quests for the wicked'
$ws.Range("C7").Value = 'A month has gone by, and you are exhausted. Emotionally, mentally, and physically. Crafting, Fighting, Enchanting, Disenchanting – all of it has worn you out.<br /> <br /> Sitting in the shade of an old oak tree in the middle of now where you bask in the warmth of the summer sun and the breeze that fills the sense with a fast-approaching autumn.<br /> <br /> “So, you have been busy child!”<br /> <br /> You turn around and see The Guide standing there. His outfit matches this time. Brown t-shirt, pants, and shoes of nondescript origin.<br /> <br /> He also looks much younger than the last time you saw him, 25 if you had to guess. But how could an 80+ year old man become so young? What magics is this?<br /> <br /> He seems thin, slight build to his frame, mostly tall and thin with a clean-cut square jaw.<br /> <br /> “What have you been up to child? Have you managed to craft yourself some better gear? Maybe enchant said gear?”<br /> <br /> You tell him all about the adventures you have been on, the Blacksmith, Enchantress and of the adventures, notoriety, and ability to create specific enchantments.<br /> <br /> “Good, I am glad you have been busy! List, I do not have all the time in the world, alas there is a rumor I must follow up on.”<br /> <br /> You begin to ask, but he raises his hand: “Listen, there is a Helpless goblin not far from here. I need you to get The Book of Guidance from him. Simple task. Think you can manage that?”<br /> <br /> You stare at him for a moment, wondering why it’s so important but none the less agree to do as he asks.'
$ws.Range("D7").Value = 'This is easy. No matter if you are on mobile or PC, click the quests tab at the top. For mobile you can scroll this section.<br /> <br /> Either way, you should see a large quest tree. Quests are how you unlock features in Tlessa. Click on Hunting Expedition quest, which is where you must start to progress down the tree.<br /> <br /> Every quest has explicit instructions on how to complete them under the requirements tab. Each requirement will have a help icon beside it, click it to reveal explicit instructions on how to get each one.<br /> <br /> On the NPC tab, is information about the NPC. If the NPC is not on surface, there will be explicit instructions on how to access the NPC.<br /> <br /> There are a lot of quests in Tlessa which is what we use to tell a bit of a story and to gate content and plane access behind.<br /> <br /> For this Guide Quest, I want you to scroll (on mobile) to the far right, on pc, look for the quest name: Goblins Accuracy.<br /> <br /> I want you to work your way down the tree to complete this quest. Easy peasy lemon squeezy!'
$ws.Range("E7").Value = 40
$ws.Range("K7").Value = 58
$ws.Range("L7").Value = 310035
$ws.Range("M7").Value = 50

# --- New row 8: "Using Quest Items" ---
$ws.Range("A8").Value = 22
$ws.Range("B8").Value = 'Using Quest Items'
$ws.Range("C8").Value = 'Walking through the streets, you flip through the book in your hands feverishly, hoping for some small semblance of what the contents were, before they magically – and before your eyes – changed.<br /> <br /> “Find the book?” Comes a familiar voice, breaking your concentration. You look up and see The Guide. You explain that you traded a Goddess Bow for this book, a book that you thought might explain something that’s been on the back of your mind. Your place in this world.<br /> <br /> You explain that the words have changed, magically, right before your eyes, as if someone wanted to hide the true meaning of their contents.<br /> <br /> The Guide stairs at you for a moment before asking for the book, which you hand over. He thumbs through the pages. “Yup, nothing but Accuracy information here. Seems useful since you still seem to be having issues hitting things. With this you should become skilled in no time.”<br /> <br /> He hands the book back to you, “You seem to have a lot of questions regarding The Creator. Why are you hunting these answers so?”<br /> <br /> You explain that since you met The Guide and other peoples on your travels, that something about this world doesn’t feel right. Something is off. You can feel a tear, a gnawing pain in the back of your neck. It’s as if the answers are all right there, but you cannot access them.<br /> <br /> “Hmmm.” The Guide begins. “I might know someone who might have some answers. As for The Creator, I have been researching his origins. Something about his “story” seems off. Why was he in Purgatory for so long and why did he have to escape? Was it a prison?”<br /> <br /> Both you and The Guide go silent for a long while, lost in thought. This isn’t a question of who am I? It’s becoming more of a question of where am I? Or is it?'
$ws.Range("D8").Value = 'This one is super easy. As you probably have seen by now, you have a lot of quest items, some you can trade with currencies to get other items that may or may not be further traded. Some effect your character, some do not.<br /> <br /> Always investigate the quest item you get to see if effects your character. The Book you got from that Goblin quest. I want you to use it.<br /> <br /> On your character sheet is a tab called Quest Items (for mobile, select Inventory Management from the drop down on the character sheet tab). Click this tab, look for the book you got, The Book of Guidance. Click it’s name to open a details modal.<br /> <br /> “With this quest item you will automatically gain 125% XP and 0% chance towards Accuracy.”<br /> <br /> I want you to train Accuracy for 2 hours (or more) until it’s level 50.<br /> <br /> Wait, how I use quest items – you scream. You don’t. They are used automatically. Train Accuracy, watch the levels go up. If you have been enchanting, enchant some Accuracy based gear for even MORE xp.'
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 50
$ws.Range("M8").Value = 55

